$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(44.91024457687308,31.1004705924814,33.44928385120062,30.50219492047529,29.4914735530908,30.58034994686383,30.5342980422547,26.39236306066022,29.304517261602,22.75133088443462,20.20505052713136,24.86361429121439,24.08433726047397,20.81723724954657,22.04450171303318,23.24122132136057,16.96383853129853,23.92276582094578,18.90531034016558,18.53584936594225,18.94988717266429,16.41575059825266,18.24137714446459,17.02956144594984,20.10821738624739,16.66087739432123,13.20646653887309,17.4868882746683,17.571693092365,20.04514663214314,16.17094603617221,13.7614884297717,13.23913302507742,11.77958678874597,14.89096318459261,15.21138665783127,14.75788296705673,13.59628116837121,15.52413503295918,15.24470016550046,14.55825393397043,11.62602657736301,13.21852259392912,11.42118202959637,9.708072204283781,11.87507064985761,11.5745201673848,12.70533541660086,11.11674282450107,13.30560760149967,11.69147923261556,10.66323855626308,13.37891643501271,12.71995897032764,11.95656595303555,10.12427590814369,11.08991099596042,11.28533242400127,11.363079853361,11.84699606070964,12.42052143102529,14.88670061579825,14.92285305213528,10.5118017807892,13.29402729500202,11.36884843133092,13.82130249812772,14.04602944071382,12.60052375596351,10.80145867158419,13.53639949653072,11.14616153407798,13.38151543737715,8.48638514080676,10.29810235062706,10.75956079182941,11.77416674345761,12.08802315026674,11.92693895276341,14.27900866815813,12.30565253705163,13.70256912576892,13.49961120949305,10.1748617830755,11.99826221713468)

$range = $ws.Range("A2:A86")
$data = New-Object 'object[,]' $newValues.Length,1
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $data[$i, 0] = $newValues[$i]
}
$range.Value = $data
